$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.541.11"
$ws.Range("E2").Value = "  -4.61%  "

$ws.Range("D3").Value = "2.962.41"
$ws.Range("E3").Value = "  -6.39%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.53%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.49%  "

$ws.Range("D9").Value = "2.972.80"
$ws.Range("E9").Value = "  -6.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.112"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.28%  "

$ws.Range("E11").Value = "  -7.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.367"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.70%  "

$ws.Range("D13").Value = "3.484.94"
$ws.Range("E13").Value = "  -6.40%  "

$ws.Range("E14").Value = "  -3.10%  "

$ws.Range("D15").Value = "61.612.59"
$ws.Range("E15").Value = "  -4.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.46%  "

$ws.Range("D17").Value = "2.973.38"
$ws.Range("E17").Value = "  -5.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000146"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "380.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.79%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.470"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.26%  "

$ws.Range("D26").Value = "3.093.40"
$ws.Range("E26").Value = "  -6.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.185"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.16%  "

$ws.Range("D29").Value = "0.0₃0930"
$ws.Range("E29").Value = "  -9.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.88%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("E32").Value = "  -5.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.34%  "

$ws.Range("E37").Value = "  -5.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -9.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.41%  "

$ws.Range("D41").Value = "2.408.62"
$ws.Range("E41").Value = "  -10.37%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.01%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.664"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0593"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0245"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0953"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.23%  "

$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "267.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.30%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.13%  "
